# Update automatico via Actualizar 02-17-2021 12-08-53
# Refresh the "Ultimo" (last-checked) timestamp column (D) for each
# availability block, simulating the roll-forward of the monitoring
# script's successive batches of timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-15: newest batch, stamped with the current run's timestamp.
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44244.50605336405
}

# Rows 16-29: becomes what used to be the newest batch (rows 2-15),
# re-stamped a hair later due to sequential execution.
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44244.48480421296
}

# Rows 30-43: becomes what used to be the middle batch (rows 16-29).
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44244.4634865625
}
